$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1229.16
$ws.Range("I19").Value = 1188
$ws.Range("J19").Value = 1281.5454
$ws.Range("K19").Value = 1188
$ws.Range("L19").Value = 1281.5454
$ws.Range("M19").Value = -1013
$ws.Range("N19").Value = -1631.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1028.6154
$ws.Range("I28").Value = 753.35
$ws.Range("J28").Value = 1946.1666
$ws.Range("K28").Value = 753.35
$ws.Range("L28").Value = 1946.1666
$ws.Range("M28").Value = -268.35
$ws.Range("N28").Value = -2916.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 7662.467
$ws.Range("I107").Value = 9342.5
$ws.Range("J107").Value = 942.3333
$ws.Range("K107").Value = 9342.5
$ws.Range("L107").Value = 942.3333
$ws.Range("M107").Value = -7422.5
$ws.Range("N107").Value = -4782.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5883247
$ws.Range("I137").Value = 846.04346
$ws.Range("K137").Value = 2538.13038
$ws.Range("M137").Value = 11.86961999999994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2200.5557
$ws.Range("I138").Value = 1131.7435
$ws.Range("J138").Value = 3937.375
$ws.Range("K138").Value = 3395.2305
$ws.Range("L138").Value = 11812.125
$ws.Range("M138").Value = 1744.7695
$ws.Range("N138").Value = -22092.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1006.0208
$ws.Range("I141").Value = 873.56525
$ws.Range("J141").Value = 4052.5
$ws.Range("K141").Value = 2620.69575
$ws.Range("L141").Value = 12157.5
$ws.Range("M141").Value = 2559.30425
$ws.Range("N141").Value = -22517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8788.441999999999
$ws.Range("I32").Value = 8084.6167
$ws.Range("J32").Value = 13011.4
$ws.Range("K32").Value = 8084.6167
$ws.Range("L32").Value = 13011.4
$ws.Range("M32").Value = -7797.6167
$ws.Range("N32").Value = -13585.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 16669740
$ws.Range("I74").Value = 26317546
$ws.Range("J74").Value = 5348
$ws.Range("K74").Value = 26317546
$ws.Range("L74").Value = 5348
$ws.Range("M74").Value = -26316672
$ws.Range("N74").Value = -7096

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 16669740
$ws.Range("I77").Value = 26317546
$ws.Range("J77").Value = 5348
$ws.Range("K77").Value = 131587730
$ws.Range("L77").Value = 26740
$ws.Range("M77").Value = -131583362
$ws.Range("N77").Value = -35476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 38755.332
$ws.Range("J95").Value = 38755.332
$ws.Range("L95").Value = 38755.332
$ws.Range("N95").Value = -44247.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8335841.5
$ws.Range("J132").Value = 3478.25
$ws.Range("L132").Value = 10434.75
$ws.Range("N132").Value = -15494.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5610.846
$ws.Range("I134").Value = 4268.8335
$ws.Range("J134").Value = 6761.143
$ws.Range("K134").Value = 12806.5005
$ws.Range("L134").Value = 20283.429
$ws.Range("M134").Value = -10271.5005
$ws.Range("N134").Value = -25353.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 46259.832
$ws.Range("J140").Value = 46259.832
$ws.Range("L140").Value = 46259.832
$ws.Range("N140").Value = -56619.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 179
$ws.Range("I51").Value = 179
$ws.Range("K51").Value = 537
$ws.Range("M51").Value = -77

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 964.9375
$ws.Range("I98").Value = 470
$ws.Range("J98").Value = 1129.9166
$ws.Range("K98").Value = 1410
$ws.Range("L98").Value = 3389.7498
$ws.Range("M98").Value = 88
$ws.Range("N98").Value = -6385.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4024.6
$ws.Range("I116").Value = 1695.6666
$ws.Range("J116").Value = 5022.7144
$ws.Range("K116").Value = 5086.9998
$ws.Range("L116").Value = 15068.1432
$ws.Range("M116").Value = -1644.9998
$ws.Range("N116").Value = -21952.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 753.8095
$ws.Range("I122").Value = 1231.6
$ws.Range("J122").Value = 319.45456
$ws.Range("K122").Value = 11084.4
$ws.Range("L122").Value = 2875.09104
$ws.Range("M122").Value = -8634.4
$ws.Range("N122").Value = -7775.09104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27388.889
$ws.Range("I70").Value = 69333.336
$ws.Range("J70").Value = 6416.6665
$ws.Range("K70").Value = 69333.336
$ws.Range("L70").Value = 6416.6665
$ws.Range("M70").Value = -69063.336
$ws.Range("N70").Value = -6956.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 27388.889
$ws.Range("I73").Value = 69333.336
$ws.Range("J73").Value = 6416.6665
$ws.Range("K73").Value = 69333.336
$ws.Range("L73").Value = 6416.6665
$ws.Range("M73").Value = -68397.336
$ws.Range("N73").Value = -8288.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 84503.5
$ws.Range("I113").Value = 101204.2
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 101204.2
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -99034.2
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3784.2144
$ws.Range("I126").Value = 2446
$ws.Range("J126").Value = 5256.25
$ws.Range("K126").Value = 7338
$ws.Range("L126").Value = 15768.75
$ws.Range("M126").Value = -4868
$ws.Range("N126").Value = -20708.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 184238.17
$ws.Range("J141").Value = 184238.17
$ws.Range("L141").Value = 184238.17
$ws.Range("N141").Value = -194598.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7414.2856
$ws.Range("I40").Value = 8828.571
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 8828.571
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -8692.571
$ws.Range("N40").Value = -6272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 26897.5
$ws.Range("J56").Value = 26530
$ws.Range("L56").Value = 26530
$ws.Range("N56").Value = -27912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 51663.332
$ws.Range("J94").Value = 51663.332
$ws.Range("L94").Value = 51663.332
$ws.Range("N94").Value = -53015.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1992.6818
$ws.Range("I100").Value = 1799.909
$ws.Range("J100").Value = 2185.4546
$ws.Range("K100").Value = 1799.909
$ws.Range("L100").Value = 2185.4546
$ws.Range("M100").Value = -1258.909
$ws.Range("N100").Value = -3267.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3734.3242
$ws.Range("I122").Value = 4033.4783
$ws.Range("J122").Value = 3242.8572
$ws.Range("K122").Value = 12100.4349
$ws.Range("L122").Value = 9728.571599999999
$ws.Range("M122").Value = -9650.4349
$ws.Range("N122").Value = -14628.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6268.1816
$ws.Range("I132").Value = 3161.75
$ws.Range("J132").Value = 12154.053
$ws.Range("K132").Value = 9485.25
$ws.Range("L132").Value = 36462.159
$ws.Range("M132").Value = -6955.25
$ws.Range("N132").Value = -41522.159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2139.4375
$ws.Range("I122").Value = 2539.25
$ws.Range("J122").Value = 940
$ws.Range("K122").Value = 7617.75
$ws.Range("L122").Value = 2820
$ws.Range("M122").Value = -5167.75
$ws.Range("N122").Value = -7720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5349.647
$ws.Range("I126").Value = 2180.3076
$ws.Range("J126").Value = 15650
$ws.Range("K126").Value = 6540.9228
$ws.Range("L126").Value = 46950
$ws.Range("M126").Value = -4070.9228
$ws.Range("N126").Value = -51890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 61428.75
$ws.Range("J141").Value = 61428.75
$ws.Range("L141").Value = 61428.75
$ws.Range("N141").Value = -71788.75
